# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (single) slide master / all slides
#   ppt/theme/theme2.xml  -> bound to the notes master
# The authored change swaps their contents: the slide master's theme becomes
# the stock "Office Theme" colour palette, while the notes master's theme
# becomes the "Integral / Red Violet" palette that used to live on the slide
# master. The font scheme and format (fill/line/effect) scheme are identical
# between the two theme parts, so the only substantive difference is the
# 12-colour colour scheme.
#
# Recolour the theme via ThemeColorScheme.Colors(i).RGB -- the supported way
# to edit theme colours through this object model. RGB values use VBA's
# RGB(r,g,b) = r + g*256 + b*65536 encoding.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# index : role      : target hex  : VBA RGB integer
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
